# Apply the numeric corrections described in the diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 8
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 1

# Row 5
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 2
$ws.Range("J5").Value = 6
$ws.Range("L5").Value = 2

# Row 6
$ws.Range("B6").Value = 0

# Row 7
$ws.Range("C7").Value = 0

# Row 8
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2

# Row 10
$ws.Range("C10").Value = 0

# Row 11
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("L11").Value = 0

# Row 13
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

# Row 14
$ws.Range("B14").Value = 3
$ws.Range("D14").Value = 0
$ws.Range("L14").Value = 1

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 6
$ws.Range("J15").Value = 5

# Row 16
$ws.Range("C16").Value = 0
$ws.Range("G16").Value = 0

# Row 17
$ws.Range("D17").Value = 0

# Row 18
$ws.Range("B18").Value = 1
$ws.Range("F18").Value = 1

# Row 20
$ws.Range("B20").Value = 0

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 5

# Row 25
$ws.Range("F25").Value = 1

# Row 29
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 2
$ws.Range("J29").Value = 5
